# Update NATMI LR-pair data (Cd200-Cd200r4) with recomputed TPM values.
# Rows 2-5 get new values, and the old rows 6-9 (MuSCs->ECs, MuSCs->Resolving-Mac,
# Resolving-Mac->ECs, Resolving-Mac->Resolving-Mac target-cluster pairs) are removed,
# since the sheet now only keeps one row per sending cluster (target cluster = Resolving-Mac).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Cd200"
$ws.Range("C2").Value = "Cd200r4"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 81.32496133333333
$ws.Range("H2").Value = 243.974884
$ws.Range("I2").Value = 0.6750769978981389
$ws.Range("J2").Value = 0.6750769978981389
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 8.848335333333333
$ws.Range("N2").Value = 26.545006
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 719.5905288477004
$ws.Range("R2").Value = 6476.314759629305
$ws.Range("S2").Value = 0.6750769978981389
$ws.Range("T2").Value = 0.6750769978981389

$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cd200"
$ws.Range("C3").Value = "Cd200r4"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.328723
$ws.Range("H3").Value = 18.986169
$ws.Range("I3").Value = 0.05253461241570551
$ws.Range("J3").Value = 0.05253461241570551
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.848335333333333
$ws.Range("N3").Value = 26.545006
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 55.99866333577933
$ws.Range("R3").Value = 503.987970022014
$ws.Range("S3").Value = 0.05253461241570551
$ws.Range("T3").Value = 0.05253461241570551

$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Cd200"
$ws.Range("C4").Value = "Cd200r4"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 32.165674
$ws.Range("H4").Value = 96.497022
$ws.Range("I4").Value = 0.2670066641690489
$ws.Range("J4").Value = 0.2670066641690489
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 8.848335333333333
$ws.Range("N4").Value = 26.545006
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 284.6126697746814
$ws.Range("R4").Value = 2561.514027972132
$ws.Range("S4").Value = 0.2670066641690489
$ws.Range("T4").Value = 0.2670066641690489

$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Cd200"
$ws.Range("C5").Value = "Cd200r4"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.648324
$ws.Range("H5").Value = 1.944972
$ws.Range("I5").Value = 0.005381725517106667
$ws.Range("J5").Value = 0.005381725517106666
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 8.848335333333333
$ws.Range("N5").Value = 26.545006
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 5.736588156648
$ws.Range("R5").Value = 51.629293409832
$ws.Range("S5").Value = 0.005381725517106667
$ws.Range("T5").Value = 0.005381725517106666


# Remove the now-unused rows 6-9 (sheet shrinks from A1:T9 to A1:T5)
$ws.Rows("6:9").Delete() | Out-Null
